# UiComponentClassDiagram.pptx - "update pptx for UI"
#
# The "Glossary" / "Window" label (two separate runs that happened to
# read back-to-back as "GlossaryWindow") is merged into a single run of
# text "GlossaryWindow".
#
# We locate the shape defensively (by its rendered text) instead of a
# hard-coded shape index, then force PowerPoint to rebuild the
# paragraph's runs by first assigning an unrelated placeholder string
# (so the multi-run -> single-run collapse actually happens) before
# writing the final text back. A direct same-text re-assignment is a
# no-op here because the concatenated text already reads
# "GlossaryWindow". Re-assigning straight to a value that still starts
# with the old text (e.g. "GlossaryWindow...") also keeps the runs
# split, so the placeholder must be unrelated to the original text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tf = $sh.TextFrame
        if ($tf.HasText) {
            if ($tf.TextRange.Text -eq "GlossaryWindow") {
                $sh.TextFrame.TextRange.Text = "#"
                $sh.TextFrame.TextRange.Text = "GlossaryWindow"
            }
        }
    }
}
